# Add PfHRP2 ML paper slide.
#
# The deck's slide 1 ("Ensemble ML for the Prediction of Artemisinin
# Resistance in Malaria") is a self-contained template group (red/white
# banner + mosquito photo + emoji + title textbox). The new slide reuses
# that exact template (same picture, same banner rectangles) but with the
# mosquito emoji and title textbox shifted and the title re-worded for a
# second paper ("Machine Learning in the Prediction of Artemisinin
# Resistance and Diagnostic Test Sensitivity in Malaria"). The cleanest
# COM equivalent of that authoring step is "duplicate slide 1, then nudge
# two shapes and retype the title" -- which is exactly what happened here
# (new slide lands right after slide 1, same shape ids as slide 1: group
# 23, rectangles 22/17, picture 15, rectangle 18, textbox 21).

$p = $ppt.ActivePresentation

$src = $p.Slides.Item(1)
$src.Duplicate() | Out-Null

$newSlide = $p.Slides.Item(2)
$grp = $newSlide.Shapes.Item(1)

# Mosquito emoji shape ("Rectangle 17" / id 18): reposition only (size
# unchanged).
$emoji = $grp.GroupItems.Item(4)
$emoji.Left = 601746 / 914400 * 72
$emoji.Top = 611233 / 914400 * 72

# Title textbox ("TextBox 20" / id 21): reposition/resize (height/top
# unchanged) and replace the four-line "Ensemble ML" title with the new
# single-line paper title.
$title = $grp.GroupItems.Item(5)
$title.Left = 2861862 / 914400 * 72
$title.Width = 4403557 / 914400 * 72
$title.TextFrame.TextRange.Text = "Machine Learning in the Prediction of Artemisinin Resistance and Diagnostic Test Sensitivity in Malaria"
